{"js": "// The paragraph \"<id>p028v_1</id>\" was split across three runs:\n//   1) \"<id>\"     - Courier New, color 7f6000, sz 18\n//   2) \"p028v_1\"  - default font, color 000000\n//   3) \"</id>\"    - Courier New, color 7f6000, sz 18\n// The edit merges them into a single run \"<id>p028v_1</id>\" that keeps\n// the formatting of the first run. Word's body.search() matches text\n// across run boundaries, returning one Range spanning all three runs;\n// replacing that range's text with itself (via insertText/\"Replace\")\n// collapses it into a single run using the leading run's formatting,\n// exactly mirroring the target OOXML.\nconst results = context.document.body.search(\"<id>p028v_1</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text '<id>p028v_1</id>' in document body\");\n}\n\nresults.items[0].insertText(\"<id>p028v_1</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p028v_1</id>\" was split across three runs:\n#   1) \"<id>\"     - Courier New, color 7f6000, sz 18\n#   2) \"p028v_1\"  - default font, color 000000\n#   3) \"</id>\"    - Courier New, color 7f6000, sz 18\n# The edit merges them into a single run \"<id>p028v_1</id>\" that keeps\n# the formatting of the first run.\n#\n# Find.Execute matches text across run boundaries, so it returns one\n# Range spanning all three runs. We keep that Range's leading slice\n# (matching the length of the original first run's text, \"<id>\")\n# untouched so its run formatting survives, delete the remainder\n# (\"p028v_1</id>\") from the document, and re-append that remainder text\n# onto the end of the first run via InsertAfter. That keeps everything\n# inside a single run and preserves the original run/paragraph\n# attributes (rsid*, xml:space) instead of Word re-creating the run\n# from scratch.\n\n$d = $word.ActiveDocument\n\n$firstRunText = \"<id>\"\n$wholeText = \"<id>p028v_1</id>\"\n\n$target = $d.Content\n$found = $target.Find.Execute($wholeText)\nif (-not $found) {\n    throw \"Could not find target text '$wholeText' in document\"\n}\n\n$restText = $wholeText.Substring($firstRunText.Length)\n\n# Range covering just the leading run's original text, e.g. \"<id>\"\n$leadRange = $target.Duplicate\n$leadRange.MoveEnd(1, $firstRunText.Length - $wholeText.Length)\n\n# Range covering everything after it, e.g. \"p028v_1</id>\"\n$tailRange = $target.Duplicate\n$tailRange.MoveStart(1, $firstRunText.Length)\n$tailRange.Delete()\n\n# Re-attach the removed text onto the first run so the three runs\n# collapse into the single merged run the target XML expects.\n$leadRange.InsertAfter($restText)\n"}
